# Insert "(de kick-off) " right after "Na het gesprek " and before
# "met de opdrachtgever" in the coach conversation report paragraph.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute(
    "Na het gesprek met de opdrachtgever",   # FindText
    $true,                                    # MatchCase
    $false,                                   # MatchWholeWord
    $false,                                   # MatchWildcards
    $false,                                   # MatchSoundsLike
    $false,                                   # MatchAllWordForms
    $true,                                    # Forward
    1,                                         # Wrap (wdFindContinue)
    $false,                                   # Format
    "Na het gesprek (de kick-off) met de opdrachtgever",  # ReplaceWith
    2                                          # Replace (wdReplaceAll)
)

if (-not $found) {
    throw "Could not find target sentence to update."
}
